$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" (column D) values are plain numeric-looking strings (e.g. "41.97")
# that Excel would otherwise auto-convert to real numbers on assignment, unlike
# the source data which stores everything as text. Force those specific cells
# to Text format first so the written value round-trips as a string, matching
# cells whose new value still contains a thousands-separator dot (e.g.
# "28.381.58") and therefore can't be parsed as a number anyway.
$textPriceCells = @(
  "D4","D5","D7","D8","D9","D11","D12","D14","D15","D16","D17","D18","D19",
  "D20","D21","D22","D24","D26","D27","D28","D29","D30","D31","D32","D33",
  "D34","D35","D36","D37","D38","D39","D40","D41","D43","D44","D45","D47",
  "D48","D49","D50","D51"
)
foreach ($addr in $textPriceCells) {
  $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "28.381.58"
$ws.Range("E2").Value = "  +1.27%  "

# Row 3
$ws.Range("D3").Value = "1.892.68"
$ws.Range("E3").Value = "  +1.64%  "

# Row 4
$ws.Range("D4").Value = "1.012"
$ws.Range("E4").Value = "  +0.87%  "

# Row 5
$ws.Range("D5").Value = "316.87"
$ws.Range("E5").Value = "  +1.49%  "

# Row 6
$ws.Range("E6").Value = "  +1.05%  "

# Row 7
$ws.Range("D7").Value = "0.5171"
$ws.Range("E7").Value = "  +1.40%  "

# Row 8
$ws.Range("D8").Value = "0.3931"
$ws.Range("E8").Value = "  +1.93%  "

# Row 9
$ws.Range("D9").Value = "0.08423"
$ws.Range("E9").Value = "  +1.38%  "

# Row 10
$ws.Range("E10").Value = "  +1.31%  "

# Row 11
$ws.Range("D11").Value = "41.97"
$ws.Range("E11").Value = "  +1.09%  "

# Row 12 / Row 13 - swapped coin order (WrappedEther <-> Polkadot)
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "6.293"
$ws.Range("E12").Value = "  +0.92%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.912.65"
$ws.Range("E13").Value = "  +2.74%  "

# Row 14
$ws.Range("D14").Value = "20.72"
$ws.Range("E14").Value = "  +0.55%  "

# Row 15
$ws.Range("D15").Value = "7.317"
$ws.Range("E15").Value = "  +1.14%  "

# Row 16
$ws.Range("D16").Value = "1.012"
$ws.Range("E16").Value = "  +0.85%  "

# Row 17
$ws.Range("D17").Value = "0.00001111"
$ws.Range("E17").Value = "  +1.24%  "

# Row 18
$ws.Range("D18").Value = "91.53"
$ws.Range("E18").Value = "  +0.69%  "

# Row 19
$ws.Range("D19").Value = "0.06749"

# Row 20
$ws.Range("D20").Value = "17.91"
$ws.Range("E20").Value = "  +1.09%  "

# Row 21
$ws.Range("D21").Value = "1.012"
$ws.Range("E21").Value = "  +0.95%  "

# Row 22
$ws.Range("D22").Value = "6.072"
$ws.Range("E22").Value = "  +0.45%  "

# Row 23
$ws.Range("D23").Value = "28.442.41"
$ws.Range("E23").Value = "  +1.39%  "

# Row 24
$ws.Range("D24").Value = "11.20"
$ws.Range("E24").Value = "  +0.71%  "

# Row 25
$ws.Range("E25").Value = "  +1.88%  "

# Row 26
$ws.Range("D26").Value = "160.21"
$ws.Range("E26").Value = "  +1.34%  "

# Row 27
$ws.Range("D27").Value = "2.484"
$ws.Range("E27").Value = "  -2.20%  "

# Row 28
$ws.Range("D28").Value = "20.78"
$ws.Range("E28").Value = "  +1.12%  "

# Row 29
$ws.Range("D29").Value = "126.10"
$ws.Range("E29").Value = "  +0.92%  "

# Row 30
$ws.Range("D30").Value = "0.1060"
$ws.Range("E30").Value = "  +0.30%  "

# Row 31
$ws.Range("D31").Value = "1.042"
$ws.Range("E31").Value = "  +0.49%  "

# Row 32
$ws.Range("D32").Value = "5.837"
$ws.Range("E32").Value = "  +0.04%  "

# Row 33
$ws.Range("D33").Value = "3.640"
$ws.Range("E33").Value = "  +1.18%  "

# Row 34
$ws.Range("D34").Value = "9.655"
$ws.Range("E34").Value = "  +2.49%  "

# Row 35
$ws.Range("D35").Value = "0.02469"
$ws.Range("E35").Value = "  +1.83%  "

# Row 36
$ws.Range("D36").Value = "0.06611"
$ws.Range("E36").Value = "  +1.11%  "

# Row 37
$ws.Range("D37").Value = "0.2220"
$ws.Range("E37").Value = "  +2.11%  "

# Row 38
$ws.Range("D38").Value = "1.206"
$ws.Range("E38").Value = "  +0.06%  "

# Row 39
$ws.Range("D39").Value = "0.6526"
$ws.Range("E39").Value = "  +0.71%  "

# Row 40
$ws.Range("D40").Value = "1.242"
$ws.Range("E40").Value = "  +1.33%  "

# Row 41
$ws.Range("D41").Value = "5.018"
$ws.Range("E41").Value = "  +0.56%  "

# Row 42
$ws.Range("E42").Value = "  +1.47%  "

# Row 43
$ws.Range("D43").Value = "0.6138"
$ws.Range("E43").Value = "  +0.35%  "

# Row 44
$ws.Range("D44").Value = "13.15"
$ws.Range("E44").Value = "  +0.57%  "

# Row 45
$ws.Range("D45").Value = "3.706"
$ws.Range("E45").Value = "  +1.29%  "

# Row 46
$ws.Range("E46").Value = "  +0.07%  "

# Row 47
$ws.Range("D47").Value = "2.026"
$ws.Range("E47").Value = "  +0.47%  "

# Row 48
$ws.Range("D48").Value = "1.241"
$ws.Range("E48").Value = "  +2.62%  "

# Row 49
$ws.Range("D49").Value = "121.65"
$ws.Range("E49").Value = "  +1.27%  "

# Row 50
$ws.Range("D50").Value = "0.06942"
$ws.Range("E50").Value = "  +1.41%  "

# Row 51
$ws.Range("D51").Value = "78.33"
$ws.Range("E51").Value = "  -0.08%  "
